# Applies the cryptos price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text so the saved value matches the source string exactly.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.226.16"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.603.76"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "212.07"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  -0.71%  "
Set-TextValue "D9" "0.0615"
$ws.Range("E9").Value = "  -0.79%  "
Set-TextValue "D10" "18.15"
$ws.Range("E10").Value = "  -1.42%  "
Set-TextValue "D11" "0.0810"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "1.822.64"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "1.594.98"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "26.208.34"
$ws.Range("E16").Value = "  +0.17%  "
Set-TextValue "D17" "61.26"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -0.08%  "
Set-TextValue "D20" "204.40"
$ws.Range("E20").Value = "  +2.90%  "
Set-TextValue "D21" "4.28"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  -2.46%  "
Set-TextValue "D23" "6.03"
$ws.Range("E23").Value = "  +0.00%  "
Set-TextValue "D24" "1.93"
$ws.Range("E24").Value = "  +11.51%  "
Set-TextValue "D25" "144.67"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E26").Value = "  -0.04%  "
Set-TextValue "D27" "0.123"
$ws.Range("E27").Value = "  -6.68%  "
Set-TextValue "D28" "15.20"
$ws.Range("E28").Value = "  -0.14%  "
Set-TextValue "D29" "6.54"
$ws.Range("E29").Value = "  +0.23%  "
Set-TextValue "D30" "0.0491"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("E31").Value = "  -0.87%  "
Set-TextValue "D32" "3.15"
$ws.Range("E32").Value = "  -0.21%  "
Set-TextValue "D33" "2.92"
$ws.Range("E33").Value = "  -4.14%  "
Set-TextValue "D34" "1.49"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "1.140.03"
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.34%  "
Set-TextValue "D40" "0.785"
$ws.Range("E40").Value = "  -0.93%  "
Set-TextValue "D41" "0.496"
$ws.Range("E41").Value = "  -2.31%  "
Set-TextValue "D42" "0.784"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "1.737.95"
$ws.Range("E44").Value = "  -0.71%  "
Set-TextValue "D45" "92.14"
$ws.Range("E45").Value = "  -1.12%  "
Set-TextValue "D46" "1.51"
$ws.Range("E46").Value = "  -2.67%  "
Set-TextValue "D47" "54.21"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "0.0₇0951"
$ws.Range("E51").Value = "  -10.54%  "
